$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-12-15T07:01:15.464963+00:00"
$ws.Range("K3").Value = "2025-12-15T07:01:15.465001+00:00"
$ws.Range("K4").Value = "2025-12-15T07:01:15.465023+00:00"
$ws.Range("K5").Value = "2025-12-15T07:01:18.220155+00:00"
$ws.Range("K6").Value = "2025-12-15T07:01:18.220185+00:00"
$ws.Range("K7").Value = "2025-12-15T07:01:18.220206+00:00"
$ws.Range("K8").Value = "2025-12-15T07:01:20.632711+00:00"
$ws.Range("K9").Value = "2025-12-15T07:01:23.375509+00:00"
$ws.Range("K10").Value = "2025-12-15T07:01:26.166509+00:00"
$ws.Range("K11").Value = "2025-12-15T07:01:28.505951+00:00"
$ws.Range("K12").Value = "2025-12-15T07:01:33.490158+00:00"
$ws.Range("K13").Value = "2025-12-15T07:01:33.490188+00:00"
$ws.Range("K14").Value = "2025-12-15T07:01:35.912137+00:00"
$ws.Range("K15").Value = "2025-12-15T07:01:38.639483+00:00"
$ws.Range("K16").Value = "2025-12-15T07:01:41.348198+00:00"
$ws.Range("K17").Value = "2025-12-15T07:01:43.695813+00:00"
$ws.Range("K18").Value = "2025-12-15T07:01:43.695843+00:00"
$ws.Range("K19").Value = "2025-12-15T07:01:43.695861+00:00"
$ws.Range("K20").Value = "2025-12-15T07:01:43.695879+00:00"
$ws.Range("K21").Value = "2025-12-15T07:01:43.695895+00:00"
$ws.Range("K22").Value = "2025-12-15T07:01:46.438713+00:00"
$ws.Range("K23").Value = "2025-12-15T07:01:46.438743+00:00"
$ws.Range("K24").Value = "2025-12-15T07:01:46.438762+00:00"
$ws.Range("K25").Value = "2025-12-15T07:01:48.648329+00:00"
$ws.Range("K26").Value = "2025-12-15T07:01:48.648363+00:00"
$ws.Range("K27").Value = "2025-12-15T07:01:48.648383+00:00"
$ws.Range("K28").Value = "2025-12-15T07:01:48.648403+00:00"
$ws.Range("K29").Value = "2025-12-15T07:01:48.648420+00:00"
$ws.Range("K30").Value = "2025-12-15T07:01:50.889440+00:00"
$ws.Range("K31").Value = "2025-12-15T07:01:50.889493+00:00"
$ws.Range("K32").Value = "2025-12-15T07:01:50.889518+00:00"
$ws.Range("K33").Value = "2025-12-15T07:01:53.683787+00:00"
$ws.Range("K34").Value = "2025-12-15T07:01:53.683817+00:00"
$ws.Range("K35").Value = "2025-12-15T07:01:53.683836+00:00"
$ws.Range("K36").Value = "2025-12-15T07:01:55.889251+00:00"
$ws.Range("K37").Value = "2025-12-15T07:01:58.178439+00:00"
$ws.Range("K38").Value = "2025-12-15T07:01:58.178469+00:00"
$ws.Range("K39").Value = "2025-12-15T07:01:58.178486+00:00"
$ws.Range("K40").Value = "2025-12-15T07:02:00.531355+00:00"
$ws.Range("K41").Value = "2025-12-15T07:02:03.323879+00:00"
$ws.Range("K42").Value = "2025-12-15T07:02:03.323913+00:00"
$ws.Range("K43").Value = "2025-12-15T07:02:06.188650+00:00"
$ws.Range("K44").Value = "2025-12-15T07:02:06.188683+00:00"
